$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 217.16667
$ws.Range("I33").Value = 146
$ws.Range("K33").Value = 146
$ws.Range("M33").Value = 83

$ws.Range("H40").Value = 1390
$ws.Range("J40").Value = 1390
$ws.Range("L40").Value = 1390
$ws.Range("N40").Value = -1740

$ws.Range("H62").Value = 6169.3076
$ws.Range("I62").Value = 6266.75
$ws.Range("K62").Value = 6266.75
$ws.Range("M62").Value = -5642.75

$ws.Range("H64").Value = 3586.279
$ws.Range("I64").Value = 2971.6785
$ws.Range("K64").Value = 2971.6785
$ws.Range("M64").Value = -2723.6785

$ws.Range("H65").Value = 6169.3076
$ws.Range("I65").Value = 6266.75
$ws.Range("K65").Value = 31333.75
$ws.Range("M65").Value = -28213.75

$ws.Range("H67").Value = 3586.279
$ws.Range("I67").Value = 2971.6785
$ws.Range("K67").Value = 2971.6785
$ws.Range("M67").Value = -2113.6785

$ws.Range("H74").Value = 10473.421
$ws.Range("I74").Value = 10473.421
$ws.Range("K74").Value = 10473.421
$ws.Range("M74").Value = -9537.421

$ws.Range("H77").Value = 10473.421
$ws.Range("I77").Value = 10473.421
$ws.Range("K77").Value = 52367.105
$ws.Range("M77").Value = -47687.105

$ws.Range("H112").Value = 2831.56
$ws.Range("J112").Value = 3284.5625
$ws.Range("L112").Value = 9853.6875
$ws.Range("N112").Value = -12069.6875

$ws.Range("H127").Value = 5058.0625
$ws.Range("I127").Value = 5174.2144
$ws.Range("K127").Value = 15522.6432
$ws.Range("M127").Value = -10562.6432

$ws.Range("H133").Value = 79999
$ws.Range("J133").Value = 79999
$ws.Range("L133").Value = 79999
$ws.Range("N133").Value = -90119

$ws.Range("H134").Value = 84499.5
$ws.Range("J134").Value = 84499.5
$ws.Range("L134").Value = 84499.5
$ws.Range("N134").Value = -94639.5

$ws.Range("H136").Value = 118910.6
$ws.Range("J136").Value = 118910.6
$ws.Range("L136").Value = 118910.6
$ws.Range("N136").Value = -129110.6

$ws.Range("H139").Value = 65099.2
$ws.Range("J139").Value = 65099.2
$ws.Range("L139").Value = 65099.2
$ws.Range("N139").Value = -75379.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1999.2307
$ws.Range("I45").Value = 2048.75
$ws.Range("K45").Value = 2048.75
$ws.Range("M45").Value = -1671.75

$ws.Range("H74").Value = 1993104.8
$ws.Range("I74").Value = 2317285.8
$ws.Range("K74").Value = 2317285.8
$ws.Range("M74").Value = -2316411.8

$ws.Range("H77").Value = 1993104.8
$ws.Range("I77").Value = 2317285.8
$ws.Range("K77").Value = 11586429
$ws.Range("M77").Value = -11582061

$ws.Range("H102").Value = 1237.3889
$ws.Range("I102").Value = 1723.3334
$ws.Range("K102").Value = 1723.3334
$ws.Range("M102").Value = -101.3334

$ws.Range("H132").Value = 643371.3
$ws.Range("J132").Value = 1641.7059
$ws.Range("L132").Value = 4925.1177
$ws.Range("N132").Value = -9985.117699999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2272.8372
$ws.Range("I105").Value = 1790.9667
$ws.Range("K105").Value = 1790.9667
$ws.Range("M105").Value = -43.96669999999995

$ws.Range("H134").Value = 5383937
$ws.Range("J134").Value = 41693164
$ws.Range("L134").Value = 125079492
$ws.Range("N134").Value = -125084562

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 16449.111
$ws.Range("I86").Value = 52000
$ws.Range("K86").Value = 52000
$ws.Range("M86").Value = -50877

$ws.Range("H89").Value = 16449.111
$ws.Range("I89").Value = 52000
$ws.Range("K89").Value = 260000
$ws.Range("M89").Value = -254384

$ws.Range("H122").Value = 21373.889
$ws.Range("I122").Value = 3455.1428
$ws.Range("K122").Value = 10365.4284
$ws.Range("M122").Value = -7915.428400000001

$ws.Range("H134").Value = 2604.093
$ws.Range("I134").Value = 2486.3823
$ws.Range("J134").Value = 3048.7778
$ws.Range("K134").Value = 7459.146900000001
$ws.Range("L134").Value = 9146.3334
$ws.Range("M134").Value = -4924.146900000001
$ws.Range("N134").Value = -14216.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 17499.5
$ws.Range("I87").Value = 8499
$ws.Range("J87").Value = 21999.75
$ws.Range("K87").Value = 25497
$ws.Range("L87").Value = 65999.25
$ws.Range("M87").Value = -24249
$ws.Range("N87").Value = -68495.25

$ws.Range("H90").Value = 17499.5
$ws.Range("I90").Value = 8499
$ws.Range("J90").Value = 21999.75
$ws.Range("K90").Value = 76491
$ws.Range("L90").Value = 197997.75
$ws.Range("M90").Value = -70251
$ws.Range("N90").Value = -210477.75

$ws.Range("H126").Value = 11891.333
$ws.Range("I126").Value = 3511
$ws.Range("K126").Value = 10533
$ws.Range("M126").Value = -5593

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 29999.5
$ws.Range("J63").Value = 29999.5
$ws.Range("L63").Value = 29999.5
$ws.Range("N63").Value = -31371.5

$ws.Range("H66").Value = 29999.5
$ws.Range("J66").Value = 29999.5
$ws.Range("L66").Value = 89998.5
$ws.Range("N66").Value = -96862.5

$ws.Range("H101").Value = 63584.5
$ws.Range("J101").Value = 63584.5
$ws.Range("L101").Value = 63584.5
$ws.Range("N101").Value = -70074.5

$ws.Range("H102").Value = 4369.5
$ws.Range("I102").Value = 4369.5
$ws.Range("K102").Value = 4369.5
$ws.Range("M102").Value = -2747.5

$ws.Range("H132").Value = 8826.755999999999
$ws.Range("I132").Value = 8212.833000000001
$ws.Range("J132").Value = 9693.471
$ws.Range("K132").Value = 24638.499
$ws.Range("L132").Value = 29080.413
$ws.Range("M132").Value = -22108.499
$ws.Range("N132").Value = -34140.413

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2536.6365
$ws.Range("I7").Value = 1779.75
$ws.Range("K7").Value = 1779.75
$ws.Range("M7").Value = -1667.75

$ws.Range("H22").Value = 1173
$ws.Range("J22").Value = 1679.2858
$ws.Range("L22").Value = 1679.2858
$ws.Range("N22").Value = -2269.2858

$ws.Range("H27").Value = 1173
$ws.Range("J27").Value = 1679.2858
$ws.Range("L27").Value = 1679.2858
$ws.Range("N27").Value = -1893.2858

$ws.Range("H93").Value = 3060.8262
$ws.Range("I93").Value = 2693.5557
$ws.Range("K93").Value = 2693.5557
$ws.Range("M93").Value = -1445.5557

$ws.Range("H101").Value = 15344.8
$ws.Range("J101").Value = 15344.8
$ws.Range("L101").Value = 15344.8
$ws.Range("N101").Value = -21834.8

$ws.Range("H122").Value = 3399.2856
$ws.Range("I122").Value = 2689.5
$ws.Range("J122").Value = 5173.75
$ws.Range("K122").Value = 8068.5
$ws.Range("L122").Value = 15521.25
$ws.Range("M122").Value = -5618.5
$ws.Range("N122").Value = -20421.25

$ws.Range("H126").Value = 2536.6365
$ws.Range("I126").Value = 1779.75
$ws.Range("K126").Value = 5339.25
$ws.Range("M126").Value = -2869.25

$ws.Range("H132").Value = 7431250
$ws.Range("I132").Value = 16715688
$ws.Range("J132").Value = 3699.8
$ws.Range("K132").Value = 50147064
$ws.Range("L132").Value = 11099.4
$ws.Range("M132").Value = -50144534
$ws.Range("N132").Value = -16159.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 14973
$ws.Range("I33").Value = 14898
$ws.Range("K33").Value = 14898
$ws.Range("M33").Value = -14648

$ws.Range("H36").Value = 14973
$ws.Range("I36").Value = 14898
$ws.Range("K36").Value = 14898
$ws.Range("M36").Value = -14648

$ws.Range("H37").Value = 23500
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 23500
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 23500
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -23906

$ws.Range("H96").Value = 37039308
$ws.Range("I96").Value = 83334740
$ws.Range("J96").Value = 2969
$ws.Range("K96").Value = 83334740
$ws.Range("L96").Value = 2969
$ws.Range("M96").Value = -83333367
$ws.Range("N96").Value = -5715

$ws.Range("H132").Value = 7939168.5
$ws.Range("I132").Value = 8335911.5
$ws.Range("K132").Value = 25007734.5
$ws.Range("M132").Value = -25005204.5
